$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 0. Remove all existing hyperlink bookkeeping first. The underlying
#    cell values/text/styles are untouched by this; only the
#    <hyperlinks> entries (and their relationships) are cleared, and
#    we will recreate them below, pointing at the post-insert cells.
# ------------------------------------------------------------------
$ws.Cells.Hyperlinks.Delete()

# ------------------------------------------------------------------
# 1. Insert a new column before column A. This shifts all existing
#    data (and column width formatting) from A/B/C to B/C/D.
# ------------------------------------------------------------------
$ws.Columns("A:A").Insert()

# ------------------------------------------------------------------
# 2. Add the new "Anal" tag cell in A29 (plain string, no hyperlink)
#    and the new video-link cell in D29.
# ------------------------------------------------------------------
$ws.Range("A29").Value = "Anal"
$ws.Range("D29").Value = "https://jav.la/video/47578/mukd-452-jav-streaming-anal-anal-sex-for-the-first-time-when-i-was-born-the-last-day-of-virginityand-the-first-anal-anal-fucking-out-yuka-shirai-yukari.html"

# ------------------------------------------------------------------
# 3. Recreate all hyperlinks at their shifted cell locations (old
#    column B -> C, old column C -> D), plus the new one on D29.
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D5"), "https://bestjavporn.com/video/pppd-812c/")
$ws.Hyperlinks.Add($ws.Range("D6"), "https://bestjavporn.com/video/pppd-779-yuria-yoshine-her-first-creampie-intimate-sex-holding-onto-her-tits-there-s-no-greater-pleasure-than-being/")
$ws.Hyperlinks.Add($ws.Range("D7"), "https://bestjavporn.com/video/ebod-707-marina-yuzuki-yuria-yoshine-m-cup-titties-and-k-cup-titties-huge-tits-a-pussy-sandwich-reverse-threesome-harlem-yuria-yoshine/")
$ws.Hyperlinks.Add($ws.Range("C9"), "https://bestjavporn.com/pornstar/marina-yuzuki/", [System.Type]::Missing, "Marina Yuzuki", "https://bestjavporn.com/pornstar/marina-yuzuki/")
$ws.Hyperlinks.Add($ws.Range("D8"), "https://bestjavporn.com/video/venu-950-streamjav-2-seconds-after-father-leaves-for-work-this-stepmother-and-stepson-will-start-fucking-yuria-yoshine/?asgtbndr=1")
$ws.Hyperlinks.Add($ws.Range("C10"), "https://bestjavporn.com/pornstar/maina-yuri/", [System.Type]::Missing, "Maina Yuri", "https://bestjavporn.com/pornstar/maina-yuri/")
$ws.Hyperlinks.Add($ws.Range("D10"), "https://bestjavporn.com/video/flav-248-free-jav-slutty-huge-ass-school-swimsuit-s-l-loves-sloppy-kisses-in-bukkake-club-maina-yuri/?asgtbndr=1")
$ws.Hyperlinks.Add($ws.Range("D14"), "https://bestjavporn.com/video/fc2-ppv-1458021/?asgtbndr=1")
$ws.Hyperlinks.Add($ws.Range("D16"), "https://bestjavporn.com/video/otim-032-javout-nene-tanaka-for-streaming-editions-daydream-pov-fantasies-she-s-plain-and-doesn-t-stand-out-in-a-crowd-but/?asgtbndr=1")
$ws.Hyperlinks.Add($ws.Range("D17"), "https://bestjavporn.com/video/jmty-032-the-bride-market-nene-not-her-real-name-nene-tanaka/")
$ws.Hyperlinks.Add($ws.Range("D18"), "https://bestjavporn.com/video/yst-224-jav-video-my-coworker-s-wife-becomes-my-personal-sex-toy-for-3-days-nene-tanaka/")
$ws.Hyperlinks.Add($ws.Range("C16"), "https://bestjavporn.com/pornstar/nene-tanaka/", [System.Type]::Missing, "Nene Tanaka", "https://bestjavporn.com/pornstar/nene-tanaka/")
$ws.Hyperlinks.Add($ws.Range("C21"), "https://bestjavporn.com/pornstar/aimi-irie/", [System.Type]::Missing, "Aimi Irie", "https://bestjavporn.com/pornstar/aimi-irie/")
$ws.Hyperlinks.Add($ws.Range("D21"), "https://bestjavporn.com/video/avop-026-lolita-stepsisters-forbidden-relations-3-secret-lesbian-acts-behind-stepmom-s-back/")
$ws.Hyperlinks.Add($ws.Range("C22"), "https://bestjavporn.com/pornstar/cocoa-aisu/", [System.Type]::Missing, "Cocoa Aisu", "https://bestjavporn.com/pornstar/cocoa-aisu/")
$ws.Hyperlinks.Add($ws.Range("D13"), "https://hpjav.tv/140943/bama-006")
$ws.Hyperlinks.Add($ws.Range("D24"), "https://bestjavporn.com/video/sdde-609-javseen-waka-misono-sumire-kurokawa-starting-today-you-are-a-member-of-the-sex-earth-protection-unit-2169-you-must-have-sex-with-the/")
$ws.Hyperlinks.Add($ws.Range("C24"), "https://bestjavporn.com/pornstar/waka-misono/", [System.Type]::Missing, "Waka Misono", "https://bestjavporn.com/pornstar/waka-misono/")
$ws.Hyperlinks.Add($ws.Range("D27"), "https://bestjavporn.com/video/jav-fitch-juny-002-hime-ichimaru/?asgtbndr=1")
$ws.Hyperlinks.Add($ws.Range("C27"), "https://bestjavporn.com/pornstar/hime-ichimaru/", [System.Type]::Missing, "Hime Ichimaru", "https://bestjavporn.com/pornstar/hime-ichimaru/")
$ws.Hyperlinks.Add($ws.Range("D29"), "https://jav.la/video/47578/mukd-452-jav-streaming-anal-anal-sex-for-the-first-time-when-i-was-born-the-last-day-of-virginityand-the-first-anal-anal-fucking-out-yuka-shirai-yukari.html")

# ------------------------------------------------------------------
# 4. Re-apply the Hyperlink style to every linked cell, since adding
#    a hyperlink programmatically can otherwise create a near-
#    duplicate style; this keeps cells on the original "Hyperlink"
#    cell style (s=1), matching the rest of the sheet.
# ------------------------------------------------------------------
foreach ($addr in @("D5","D6","D7","C9","D8","C10","D10","D14","D16","D17","D18","C16","C21","D21","C22","D13","D24","C24","D27","C27","D29")) {
    $ws.Range($addr).Style = "Hyperlink"
}

# ------------------------------------------------------------------
# 5. Update the active selection to match the final workbook state.
# ------------------------------------------------------------------
$ws.Range("D29").Select()
